# Updated capital structure database
# Applies the data refresh for the two Saudi Arabia / Machinery rows (rows 2 & 3)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2,3) {
    $ws.Range("D$r").Value = -0.278

    $ws.Range("G$r").Value = -0.1437206178643385
    $ws.Range("H$r").Value = -0.1437206178643385
    $ws.Range("I$r").Value = 0.06178643384822027
    $ws.Range("J$r").Value = 0.06178643384822027
    $ws.Range("K$r").Value = -75
    $ws.Range("L$r").Value = -0.5036937541974479

    $ws.Range("U$r").Value = 19
    $ws.Range("V$r").Value = 0.1760889712696941
    $ws.Range("W$r").Value = -0.8426966292134831
    $ws.Range("X$r").Value = 0.2041449860151529
    $ws.Range("Y$r").Value = -1.046841615228636
    $ws.Range("Z$r").Value = 0.3786876907426246
    $ws.Range("AA$r").Value = 0.02339776195320447
    $ws.Range("AB$r").Value = 0.07080228571312558
    $ws.Range("AC$r").Value = -0.0474045237599211
    $ws.Range("AD$r").Value = 318.4
    $ws.Range("AE$r").Value = 0
    $ws.Range("AF$r").Value = 318.4
    $ws.Range("AG$r").Value = 299.4
    $ws.Range("AH$r").Value = 0.7468918601923529
    $ws.Range("AI$r").Value = 0.8844444444444444
    $ws.Range("AJ$r").Value = 0.7350847041492757
    $ws.Range("AK$r").Value = 0.8780058651026392
    $ws.Range("AL$r").Value = 0
    $ws.Range("AM$r").Value = -2.17
    $ws.Range("AN$r").Value = -20.81045751633987

    # Column AO (ebit_interest_expenses) no longer has a value for this row
    $ws.Range("AO$r").ClearContents()

    $ws.Range("AP$r").Value = -19.56862745098039
    $ws.Range("AQ$r").Value = -4.239631336405529
}
